$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.398.80'
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').Value = '2.532.04'
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.88'
$ws.Range('E5').Value = '  -1.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.47'
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.566'
$ws.Range('E7').Value = '  -1.58%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.527'
$ws.Range('E9').Value = '  -2.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.66'
$ws.Range('E10').Value = '  +0.31%  '
$ws.Range('E11').Value = '  -1.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.32'
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').Value = '2.923.02'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.68'
$ws.Range('E15').Value = '  +4.38%  '
$ws.Range('D16').Value = '2.555.41'
$ws.Range('E16').Value = '  -3.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.827'
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('D18').Value = '42.436.66'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.78'
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('D20').Value = '0.0₃0948'
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.18'
$ws.Range('E21').Value = '  -3.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.05'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '243.63'
$ws.Range('E23').Value = '  -3.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.89'
$ws.Range('E24').Value = '  -2.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.04'
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.89'
$ws.Range('E27').Value = '  -3.34%  '
$ws.Range('E28').Value = '  -4.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.17'
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.11'
$ws.Range('E30').Value = '  -0.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.51'
$ws.Range('E31').Value = '  +2.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.68'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('E33').Value = '  +15.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0794'
$ws.Range('E34').Value = '  -1.62%  '
$ws.Range('E35').Value = '  -3.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.01'
$ws.Range('E36').Value = '  -5.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.15'
$ws.Range('E37').Value = '  -8.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.05'
$ws.Range('E38').Value = '  -5.13%  '
$ws.Range('E39').Value = '  -0.63%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.26'
$ws.Range('E41').Value = '  +9.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.67'
$ws.Range('E42').Value = '  -3.99%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('E44').Value = '  +0.72%  '
$ws.Range('E45').Value = '  -2.58%  '
$ws.Range('D46').Value = '1.961.33'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.92'
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('D48').Value = '2.777.03'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '80.70'
$ws.Range('E49').Value = '  -3.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.191'
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.843'
$ws.Range('E51').Value = '  +7.62%  '
